$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.838.96"
$ws.Range("E2").Value = "  -4.08%  "

$ws.Range("D3").Value = "2.445.34"
$ws.Range("E3").Value = "  -3.24%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.75"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.88"
$ws.Range("E6").Value = "  -7.93%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.547"
$ws.Range("E7").Value = "  -3.28%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  -6.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.78"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0767"
$ws.Range("E11").Value = "  -4.22%  "

$ws.Range("E12").Value = "  -1.10%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.826.52"
$ws.Range("E13").Value = "  -3.08%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.86"
$ws.Range("E14").Value = "  -6.72%  "

$ws.Range("D15").Value = "2.476.28"
$ws.Range("E15").Value = "  -2.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.60"
$ws.Range("E16").Value = "  -4.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.768"
$ws.Range("E17").Value = "  -4.83%  "

$ws.Range("D18").Value = "40.870.54"
$ws.Range("E18").Value = "  -4.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.19"
$ws.Range("E19").Value = "  -7.55%  "

$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  -4.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.93"
$ws.Range("E21").Value = "  -9.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.29"
$ws.Range("E22").Value = "  -3.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.55"
$ws.Range("E23").Value = "  -4.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.72"
$ws.Range("E24").Value = "  -5.21%  "

$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("E26").Value = "  -7.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.45"
$ws.Range("E27").Value = "  -7.83%  "

$ws.Range("E28").Value = "  -6.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.45"
$ws.Range("E29").Value = "  -6.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.25"
$ws.Range("E30").Value = "  -7.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "150.16"
$ws.Range("E31").Value = "  -4.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.38"
$ws.Range("E32").Value = "  -5.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.65"
$ws.Range("E33").Value = "  -5.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.53"
$ws.Range("E34").Value = "  -3.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0728"
$ws.Range("E35").Value = "  -6.81%  "

$ws.Range("E36").Value = "  -6.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.53"
$ws.Range("E37").Value = "  -7.65%  "

$ws.Range("E38").Value = "  -6.97%  "

$ws.Range("E39").Value = "  -3.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.100"
$ws.Range("E40").Value = "  -9.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.08"
$ws.Range("E41").Value = "  -3.12%  "

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.56"
$ws.Range("E43").Value = "  -12.27%  "

$ws.Range("D44").Value = "1.951.95"
$ws.Range("E44").Value = "  -2.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0279"
$ws.Range("E45").Value = "  -6.72%  "

$ws.Range("E46").Value = "  -8.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.49"
$ws.Range("E47").Value = "  -4.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "68.21"
$ws.Range("E48").Value = "  -5.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "95.20"
$ws.Range("E49").Value = "  -5.35%  "

$ws.Range("E50").Value = "  -7.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.86"
$ws.Range("E51").Value = "  -7.85%  "
